# Update the cosinor analysis results on Sheet1 (row 2 and row 3)
# to reflect the re-run CircadiPy simulation analysis values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("G2").Value = [double]"2.642330798607873e-14"
$ws.Range("H2").Value = [double]"9.234235950850508e-14"
$ws.Range("K2").Value = [double]"41.4391415716992"
$ws.Range("L2").Value = "[30.83457479692757, 52.04370834647082]"
$ws.Range("M2").Value = [double]"7.043254868221993e-13"
$ws.Range("N2").Value = [double]"1.408650973644399e-12"
$ws.Range("O2").Value = [double]"1.641552918091963"
$ws.Range("P2").Value = "[1.339658128557808, 1.9434477076261185]"
$ws.Range("S2").Value = [double]"59.91288495659356"
$ws.Range("T2").Value = "[53.126795129632654, 66.69897478355446]"
$ws.Range("W2").Value = [double]"18.39459459459493"
$ws.Range("X2").Value = [double]"17.19819819819851"
$ws.Range("Y2").Value = [double]"19.59099099099135"

# --- Row 3 ---
$ws.Range("E3").Value = [double]"24.17000000000034"
$ws.Range("G3").Value = [double]"6.661338147750939e-16"
$ws.Range("H3").Value = [double]"1.125859968633962e-14"
$ws.Range("K3").Value = [double]"45.62360837333575"
$ws.Range("L3").Value = "[32.87865825141181, 58.36855849525969]"
$ws.Range("M3").Value = [double]"3.059508202341021e-11"
$ws.Range("N3").Value = [double]"3.059508202341021e-11"
$ws.Range("O3").Value = [double]"0.4842895582110396"
$ws.Range("P3").Value = "[0.1823947686768861, 0.7861843477451931]"
$ws.Range("Q3").Value = [double]"0.001811525057478391"
$ws.Range("R3").Value = [double]"0.001811525057478391"
$ws.Range("S3").Value = [double]"60.89204763652153"
$ws.Range("T3").Value = "[54.00262487998698, 67.78147039305608]"
$ws.Range("W3").Value = [double]"22.30704704704736"
$ws.Range("X3").Value = [double]"21.14572572572602"
$ws.Range("Y3").Value = [double]"23.46836836836869"
